$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 - RGL1
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
# bump the sheet's stored max row-outline level (1 -> 2) without leaving a
# visible outline group behind: group/ungroup a throwaway row below the
# real data, then delete it.
$ws.Rows.Item(10).OutlineLevel = 2
$ws.Rows.Item(10).Delete()
$ws.Cells.Item(3,7).Value = 1
$ws.Cells.Item(3,8).Value = 2
$ws.Cells.Item(3,9).Value = 3
$ws.Cells.Item(3,10).Value = 4
$ws.Cells.Item(3,11).Value = 5
$ws.Cells.Item(3,12).Value = 6
$ws.Cells.Item(3,13).Value = 7
$ws.Range("N3").Select()

# ---------------------------------------------------------------------------
# Sheet 2 - RGL2 (becomes the active tab)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).Clear()
$ws.Cells.Item(3,7).Value = 1
$ws.Cells.Item(3,8).Value = 2
$ws.Cells.Item(3,9).Value = 3
$ws.Cells.Item(3,10).Value = 4
$ws.Cells.Item(3,11).Value = 5
$ws.Cells.Item(3,12).Value = 6
$ws.Range("N9").Select()
$ws.Activate()

# ---------------------------------------------------------------------------
# Sheet 3 - SLIM (loses the active tab)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(3)
$ws.Rows.Item(10).OutlineLevel = 2
$ws.Rows.Item(10).Delete()
$ws.Cells.Item(3,7).Value = 1
$ws.Cells.Item(3,8).Value = 2
$ws.Cells.Item(3,9).Value = 3
$ws.Range("J3").Select()

# ---------------------------------------------------------------------------
# Sheet 4 - RGLJ
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Rows.Item(10).OutlineLevel = 2
$ws.Rows.Item(10).Delete()
$ws.Cells.Item(3,7).Value = 1
$ws.Cells.Item(3,8).Value = 2
$ws.Cells.Item(3,9).Value = 3
$ws.Range("I3").Select()

# ---------------------------------------------------------------------------
# Sheet 5 - D.PANT
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)
$ws.Rows.Item(10).OutlineLevel = 3
$ws.Rows.Item(10).Delete()
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = 2
$ws.Cells.Item(4,9).Value = 3
$ws.Range("I4").Select()

# ---------------------------------------------------------------------------
# Sheet 6 - C.PANT
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(6)
$ws.Rows.Item(10).OutlineLevel = 2
$ws.Rows.Item(10).Delete()
$ws.Cells.Item(3,7).Value = 1
$ws.Cells.Item(3,8).Value = 2
$ws.Cells.Item(3,9).Value = 3
$ws.Range("A1").Select()
$ws.Range("I3").Select()

# ---------------------------------------------------------------------------
# Sheet 7 - Sheet1
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(7)
$ws.Range("H9").Select()

# ---------------------------------------------------------------------------
# workbook-level: RGL2 (tab index 1, zero-based) is now the active sheet
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(2).Activate()
